# B6-PowerPoint.pptx edit — Thu, Jul 09, 2020 3:07:37 PM
#
# 1) Re-style the three summary tables (balance-sheet slides) with the
#    built-in "no style / table grid" table style instead of the custom
#    Google-Slides-imported style.
# 2) Swap the deck's colour theme from the old "Integral" (Red Violet)
#    palette to the standard Office palette.

$p = $ppt.ActivePresentation

# --- 1. Tables: apply the new table style GUID -----------------------
$newStyleId = "{564308AF-1F1B-47D6-8DBE-04F1F3317FAF}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Theme: switch the colour scheme to the Office palette --------
$master = $p.Slides.Item(1).Master
$colors = $master.Theme.ThemeColorScheme

$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
